$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a brand-new "item" row above the totals row (row 19) so the totals
# row and the footer row shift down by one (19->20, 20->21), mirroring the
# previous item rows (7..18) that already exist on the sheet.
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).Insert()

# Clone the formatting (styles) of the last existing item row (18) into the
# freshly inserted row 19 so borders/fonts/fills/number-formats line up with
# the rest of the item rows.
$ws.Range("A18:Q18").Copy()
$ws.Range("A19:Q19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Helper used below: write a string value into a cell while preserving the
# cell's existing NumberFormat (temporarily switching to text "@" so Excel
# stores the value as a genuine string instead of re-interpreting a
# numeric-looking string such as "15.0000" or "4:0" as a number/date).
function Set-TextValue($rng, $text) {
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = $fmt
}

# New item row values (item #13 - "ليفه")
$ws.Range("A19").Value = 13
Set-TextValue $ws.Range("C19") "ليفه"
Set-TextValue $ws.Range("H19") "4:0"
Set-TextValue $ws.Range("L19") "0"
Set-TextValue $ws.Range("N19") "15.00"
Set-TextValue $ws.Range("P19") "15.0000"
Set-TextValue $ws.Range("Q19") "1:0"

# Re-create the merged cells for the new item row, matching the pattern used
# by every other item row (7..18).
$ws.Range("A19:B19").Merge()
$ws.Range("C19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("N19:O19").Merge()

# Row heights: the new item row gets the standard item-row height, and the
# totals row (now row 20) switches to the shorter height seen in the diff.
$ws.Rows.Item(19).RowHeight = 25.5
$ws.Rows.Item(20).RowHeight = 24.75

# Update the grand total to include the new item's price (676.43 + 15.00).
$ws.Range("P20").Value = 691.42999999999995

# Update the generated timestamp in the footer (now row 21) to reflect the
# new save time.
$ws.Range("A21").Value = "Monday, 14 July, 2025 10:43 AM"
